$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2006802721088435
$ws.Range("C2").Value = 0.5476190476190477
$ws.Range("J2").Value = 0.01020408163265306
$ws.Range("O2").Value = 0.003401360544217687
$ws.Range("P2").Value = 0.1564625850340136
$ws.Range("S2").Value = 0.08163265306122448
$ws.Range("B3").Value = 0.005882352941176471
$ws.Range("C3").Value = 0.05294117647058823
$ws.Range("J3").Value = 0.04117647058823529
$ws.Range("P3").Value = 0.6941176470588235
$ws.Range("S3").Value = 0.2058823529411765
$ws.Range("J4").Value = 0.02439024390243903
$ws.Range("O4").Value = 0.02439024390243903
$ws.Range("P4").Value = 0.7317073170731707
$ws.Range("S4").Value = 0.2195121951219512
$ws.Range("B6").Value = 0.06132075471698113
$ws.Range("D6").Value = 0.009433962264150943
$ws.Range("F6").Value = 0.0660377358490566
$ws.Range("J6").Value = 0.1886792452830189
$ws.Range("O6").Value = 0.04245283018867924
$ws.Range("Q6").Value = 0.1650943396226415
$ws.Range("R6").Value = 0.1132075471698113
$ws.Range("S6").Value = 0.3537735849056604
$ws.Range("B7").Value = 0.09734513274336283
$ws.Range("D7").Value = 0.01769911504424779
$ws.Range("F7").Value = 0.03097345132743363
$ws.Range("J7").Value = 0.1858407079646018
$ws.Range("O7").Value = 0.01327433628318584
$ws.Range("Q7").Value = 0.1371681415929203
$ws.Range("R7").Value = 0.08849557522123894
$ws.Range("S7").Value = 0.4292035398230089
$ws.Range("B8").Value = 0.08521739130434783
$ws.Range("D8").Value = 0.01913043478260869
$ws.Range("F8").Value = 0.07478260869565218
$ws.Range("J8").Value = 0.1008695652173913
$ws.Range("O8").Value = 0.01913043478260869
$ws.Range("Q8").Value = 0.1860869565217391
$ws.Range("R8").Value = 0.08521739130434783
$ws.Range("S8").Value = 0.4295652173913043
$ws.Range("B9").Value = 0.1278195488721804
$ws.Range("D9").Value = 0.04511278195488722
$ws.Range("F9").Value = 0.04511278195488722
$ws.Range("J9").Value = 0.1578947368421053
$ws.Range("O9").Value = 0.02255639097744361
$ws.Range("Q9").Value = 0.1804511278195489
$ws.Range("R9").Value = 0.07518796992481203
$ws.Range("S9").Value = 0.3458646616541353
$ws.Range("B10").Value = 0.102861562258314
$ws.Range("D10").Value = 0.01624129930394431
$ws.Range("F10").Value = 0.05259087393658159
$ws.Range("J10").Value = 0.1291569992266048
$ws.Range("O10").Value = 0.01778808971384378
$ws.Range("Q10").Value = 0.234338747099768
$ws.Range("R10").Value = 0.07965970610982212
$ws.Range("S10").Value = 0.3673627223511214
$ws.Range("F11").Value = 0.002857142857142857
$ws.Range("G11").Value = 0.1685714285714286
$ws.Range("J11").Value = 0.08
$ws.Range("K11").Value = 0.1914285714285714
$ws.Range("L11").Value = 0.54
$ws.Range("S11").Value = 0.01714285714285714
$ws.Range("G12").Value = 0.7397959183673469
$ws.Range("J12").Value = 0.1683673469387755
$ws.Range("K12").Value = 0.02040816326530612
$ws.Range("L12").Value = 0.03571428571428571
$ws.Range("S12").Value = 0.03571428571428571
$ws.Range("G13").Value = 0.5652173913043478
$ws.Range("J13").Value = 0.391304347826087
$ws.Range("S13").Value = 0.04347826086956522
$ws.Range("F15").Value = 0.03846153846153846
$ws.Range("H15").Value = 0.1752136752136752
$ws.Range("I15").Value = 0.02136752136752137
$ws.Range("J15").Value = 0.3034188034188034
$ws.Range("K15").Value = 0.06837606837606838
$ws.Range("M15").Value = 0.02564102564102564
$ws.Range("O15").Value = 0.09401709401709402
$ws.Range("S15").Value = 0.2735042735042735
$ws.Range("F16").Value = 0.01030927835051546
$ws.Range("H16").Value = 0.211340206185567
$ws.Range("I16").Value = 0.04639175257731959
$ws.Range("J16").Value = 0.3969072164948453
$ws.Range("K16").Value = 0.1082474226804124
$ws.Range("M16").Value = 0.02577319587628866
$ws.Range("N16").Value = 0.005154639175257732
$ws.Range("O16").Value = 0.04639175257731959
$ws.Range("S16").Value = 0.1494845360824742
$ws.Range("F17").Value = 0.01996007984031936
$ws.Range("H17").Value = 0.2115768463073852
$ws.Range("I17").Value = 0.08183632734530938
$ws.Range("J17").Value = 0.4131736526946108
$ws.Range("K17").Value = 0.09580838323353294
$ws.Range("M17").Value = 0.01596806387225549
$ws.Range("O17").Value = 0.07385229540918163
$ws.Range("S17").Value = 0.08782435129740519
$ws.Range("F18").Value = 0.01470588235294118
$ws.Range("H18").Value = 0.2107843137254902
$ws.Range("I18").Value = 0.03431372549019608
$ws.Range("J18").Value = 0.3872549019607843
$ws.Range("K18").Value = 0.1176470588235294
$ws.Range("M18").Value = 0.0392156862745098
$ws.Range("O18").Value = 0.07843137254901961
$ws.Range("S18").Value = 0.1176470588235294
$ws.Range("F19").Value = 0.01534919416730622
$ws.Range("H19").Value = 0.2640061396776669
$ws.Range("I19").Value = 0.05525709900230238
$ws.Range("J19").Value = 0.3591711435149655
$ws.Range("K19").Value = 0.1273983115886416
$ws.Range("M19").Value = 0.01841903300076746
$ws.Range("N19").Value = 0.001534919416730622
$ws.Range("O19").Value = 0.056792018419033
$ws.Range("S19").Value = 0.1020721412125863
